$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.223.48"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "'1.659.90"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "'219.08"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'0.5222"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "'0.2666"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.06334"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "'21.25"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'0.07763"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'4.443"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "'1.654.28"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "'0.5493"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "'0.0₅8271"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "'65.05"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").Value = "'26.243.02"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'4.699"
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("D20").Value = "'193.35"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").Value = "'10.21"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").Value = "'6.146"
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("D23").Value = "'1.006"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'138.83"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").Value = "'0.1240"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'7.291"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").Value = "'16.19"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "'1.415"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").Value = "'0.06069"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "'1.286"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").Value = "'3.565"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'3.359"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").Value = "'1.656"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("D34").Value = "'0.9867"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "'2.781"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.410"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").Value = "'0.5964"
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("D38").Value = "'0.01600"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "'5.981"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "'0.8661"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "'1.050.06"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "'1.003"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "'100.07"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "'1.796.15"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "'0.0₈110"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").Value = "'57.46"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").Value = "'8.136"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.482"
$ws.Range("E49").Value = "  +4.07%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05183"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").Value = "'0.4232"
$ws.Range("E51").Value = "  +0.26%  "
